$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.505.55"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.827.85"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'316.54"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "'0.5171"
$ws.Range("E7").Value = "  +2.00%  "
$ws.Range("D8").Value = "'0.3879"
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("D9").Value = "'0.08294"
$ws.Range("E9").Value = "  +7.74%  "
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("D11").Value = "'41.91"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").Value = "'6.404"
$ws.Range("E12").Value = "  +2.07%  "
$ws.Range("D13").Value = "'21.24"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").Value = "'1.004"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D16").Value = "1.829.01"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").Value = "'94.13"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("E18").Value = "  +3.63%  "
$ws.Range("D19").Value = "'0.06650"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").Value = "'17.85"
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("E22").Value = "  -1.57%  "
$ws.Range("D23").Value = "28.550.40"
$ws.Range("D24").Value = "'11.46"
$ws.Range("E24").Value = "  +3.08%  "
$ws.Range("D25").Value = "'2.255"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("E26").Value = "  +2.64%  "
$ws.Range("D27").Value = "'159.32"
$ws.Range("E27").Value = "  +1.42%  "
$ws.Range("D28").Value = "2.037.78"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").Value = "'2.427"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").Value = "'126.17"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("D32").Value = "'1.102"
$ws.Range("E32").Value = "  -2.74%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'5.744"
$ws.Range("E33").Value = "  +1.19%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.07582"
$ws.Range("E34").Value = "  +6.98%  "
$ws.Range("D35").Value = "'3.684"
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("D36").Value = "'0.2235"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("D37").Value = "'0.02378"
$ws.Range("E37").Value = "  +2.20%  "
$ws.Range("D38").Value = "'5.295"
$ws.Range("E38").Value = "  +2.71%  "
$ws.Range("D39").Value = "'12.01"
$ws.Range("E39").Value = "  +7.26%  "
$ws.Range("D40").Value = "'8.788"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("D41").Value = "'0.6400"
$ws.Range("E41").Value = "  +2.31%  "
$ws.Range("D42").Value = "'1.194"
$ws.Range("E42").Value = "  +0.69%  "
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").Value = "'13.74"
$ws.Range("E44").Value = "  +2.43%  "
$ws.Range("D45").Value = "'0.6165"
$ws.Range("E45").Value = "  +4.26%  "
$ws.Range("D46").Value = "'3.802"
$ws.Range("E46").Value = "  +2.26%  "
$ws.Range("D47").Value = "'128.02"
$ws.Range("E47").Value = "  +2.48%  "
$ws.Range("D48").Value = "'2.003"
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("D49").Value = "'1.205"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").Value = "'0.06982"
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("E51").Value = "  +0.43%  "
